# Update Generic Excel DataLayer (iPasXL)
# Updates the data table on the "Equipment" sheet with new sample values
# and resets the active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Equipment")

# Table headers (row 1) are unchanged:
# A: Tag  B: Description  C: PumpType  D: PumpDriverType
# E: DesignTemp  F: DesignPressure  G: Capacity  H: SpecificGravity  I: DifferentialPressure

# Row 2 - Equip-001
$ws.Range("B2").Value = "DESC-9"
$ws.Range("C2").Value = "PT-8"
$ws.Range("D2").Value = "PDT-3"
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 6
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 2

# Row 3 - Equip-002
$ws.Range("B3").Value = "DESC-3"
$ws.Range("C3").Value = "PT-7"
$ws.Range("D3").Value = "PDT-9"
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 4
$ws.Range("H3").Value = 7
$ws.Range("I3").Value = 9

# Row 4 - Equip-003
$ws.Range("B4").Value = "DESC-2"
$ws.Range("C4").Value = "PT-7"
$ws.Range("D4").Value = "PDT-3"
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 8

# Row 5 - Equip-004
$ws.Range("B5").Value = "DESC-3"
$ws.Range("C5").Value = "PT-3"
$ws.Range("D5").Value = "PDT-8"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 3

# Reset the active selection to A2
$ws.Range("A2").Select()
